$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 42619.892627314817
$ws.Range("A14").NumberFormat = "m/d/yy h:mm"

$ws.Range("B14").Value = 32
$ws.Range("C14").Value = 64
$ws.Range("D14").Value = 32
$ws.Range("E14").Value = 64
$ws.Range("F14").Value = 28
$ws.Range("G14").Value = 17589
$ws.Range("H14").Value = 17092
$ws.Range("I14").Value = 2818
$ws.Range("J14").Value = 427
$ws.Range("K14").Value = 217
$ws.Range("L14").Value = 39
$ws.Range("M14").Value = 16
$ws.Range("N14").Value = "Noun"
